$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date serials for A2:A33 (replacing the old ones)
$dates = @(43781,43782,43783,43784,43787,43788,43789,43790,43791,43794,43795,43796,43797,43798,43801,43802,43803,43804,43805,43808,43809,43810,43811,43812,43815,43816,43817,43818,43819,43822,43825,43826)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = 2550
    $ws.Cells.Item($row, 3).Value = 0
}
